$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 119; this shifts the existing rows 119-122
# down to 120-123, carrying their content and formatting along.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new weekly record.
$ws.Cells.Item(119, 1).Value = 7
$ws.Cells.Item(119, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(119, 3).Value = "Ñuble"
$ws.Cells.Item(119, 4).Value = 44448
$ws.Cells.Item(119, 5).Value = 16
$ws.Cells.Item(119, 6).Value = 100112032
$ws.Cells.Item(119, 7).Value = "Zapallo italiano"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 160
$ws.Cells.Item(119, 11).Value = 16000
$ws.Cells.Item(119, 12).Value = 17000
$ws.Cells.Item(119, 13).Value = 16500
$ws.Cells.Item(119, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(119, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(119, 16).Value = 330
$ws.Cells.Item(119, 17).Value = 50
$ws.Cells.Item(119, 18).Value = "Hortaliza"
